$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 3 rows (old self-cluster pair rows 8-10), shrinking the table to 6 data rows
$ws.Rows.Item(8).Resize(3).Delete()

# Update remaining data rows (2-7) with refreshed TPM-derived values
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.016376
$ws.Range("N2").Value = 0.049128
$ws.Range("O2").Value = 0.4917717717717718
$ws.Range("P2").Value = 0.4917717717717718
$ws.Range("Q2").Value = 1.359006264056
$ws.Range("R2").Value = 12.231056376504
$ws.Range("S2").Value = 0.2207811393220665
$ws.Range("T2").Value = 0.2207811393220665
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.016924
$ws.Range("N3").Value = 0.050772
$ws.Range("O3").Value = 0.5082282282282282
$ws.Range("P3").Value = 0.5082282282282282
$ws.Range("Q3").Value = 1.404483513244
$ws.Range("R3").Value = 12.640351619196
$ws.Range("S3").Value = 0.2281692722207287
$ws.Range("T3").Value = 0.2281692722207287
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("H4").Value = 189.421768
$ws.Range("I4").Value = 0.3415807409566563
$ws.Range("J4").Value = 0.3415807409566563
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.016376
$ws.Range("N4").Value = 0.049128
$ws.Range("O4").Value = 0.4917717717717718
$ws.Range("P4").Value = 0.4917717717717718
$ws.Range("Q4").Value = 1.033990290922667
$ws.Range("R4").Value = 9.305912618303999
$ws.Range("S4").Value = 0.1679797661833695
$ws.Range("T4").Value = 0.1679797661833695
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("H5").Value = 189.421768
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.016924
$ws.Range("N5").Value = 0.050772
$ws.Range("O5").Value = 0.5082282282282282
$ws.Range("P5").Value = 0.5082282282282282
$ws.Range("Q5").Value = 1.068591333877333
$ws.Range("R5").Value = 9.617322004895998
$ws.Range("S5").Value = 0.1736009747732868
$ws.Range("T5").Value = 0.1736009747732868
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 38.719942
$ws.Range("H6").Value = 116.159826
$ws.Range("I6").Value = 0.2094688475005485
$ws.Range("J6").Value = 0.2094688475005485
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.016376
$ws.Range("N6").Value = 0.049128
$ws.Range("O6").Value = 0.4917717717717718
$ws.Range("P6").Value = 0.4917717717717718
$ws.Range("Q6").Value = 0.634077770192
$ws.Range("R6").Value = 5.706699931728
$ws.Range("S6").Value = 0.1030108662663358
$ws.Range("T6").Value = 0.1030108662663358
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 38.719942
$ws.Range("H7").Value = 116.159826
$ws.Range("I7").Value = 0.2094688475005485
$ws.Range("J7").Value = 0.2094688475005485
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.016924
$ws.Range("N7").Value = 0.050772
$ws.Range("O7").Value = 0.5082282282282282
$ws.Range("P7").Value = 0.5082282282282282
$ws.Range("Q7").Value = 0.655296298408
$ws.Range("R7").Value = 5.897666685672
$ws.Range("S7").Value = 0.1064579812342127
$ws.Range("T7").Value = 0.1064579812342127
